# ---------------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# The workbook gains a new "2022-Q3" worksheet (fund-holdings detail, placed
# between "总计" and "2022-Q2"), and the "总计" summary sheet gets a new
# summary row for 2022-Q3 (the former 2022-Q2 summary row is pushed down).
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2    = $wb.Worksheets.Item("2022-Q2")

# --- 1. Create the new "2022-Q3" sheet, positioned right before "2022-Q2" ---
#     (final tab order: 总计, 2022-Q3, 2022-Q2)
$wsQ3 = $wb.Worksheets.Add($wsQ2)
$wsQ3.Name = "2022-Q3"

# --- 2. Header row (B1:H1) ---
$headersQ3 = New-Object 'object[,]' 1,7
$headersQ3[0,0] = "基金代码"
$headersQ3[0,1] = "基金名称"
$headersQ3[0,2] = "基金规模"
$headersQ3[0,3] = "股票总仓位"
$headersQ3[0,4] = "仓位占比"
$headersQ3[0,5] = "持有市值(亿元)"
$headersQ3[0,6] = "仓位排名"
$wsQ3.Range("B1:H1").Value = $headersQ3

# --- 3. Fund detail rows (B2:H8) ---
#     Columns B-G are stored as text (matches the source workbook's
#     convention of keeping numeric-looking figures like "12.06" as text);
#     column H (ranking) is numeric. A leading "'" forces text entry without
#     touching cell formatting/number format.
    $dataQ3 = New-Object 'object[,]' 7,7
    $dataQ3[0,0] = "'001668"
    $dataQ3[0,1] = "汇添富全球移动互联灵活配置混合（QDII）A"
    $dataQ3[0,2] = "'12.06"
    $dataQ3[0,3] = "'90.88"
    $dataQ3[0,4] = "'2.86"
    $dataQ3[0,5] = "'0.3449"
    $dataQ3[0,6] = 8
    $dataQ3[1,0] = "'161128"
    $dataQ3[1,1] = "易方达标普信息科技指数（QDII-LOF）人民币"
    $dataQ3[1,2] = "'4.99"
    $dataQ3[1,3] = "'91.96"
    $dataQ3[1,4] = "'1.66"
    $dataQ3[1,5] = "'0.0828"
    $dataQ3[1,6] = 9
    $dataQ3[2,0] = "'012868"
    $dataQ3[2,1] = "易方达标普信息科技指数（QDII-LOF）人民币 C"
    $dataQ3[2,2] = "'4.99"
    $dataQ3[2,3] = "'91.96"
    $dataQ3[2,4] = "'1.66"
    $dataQ3[2,5] = "'0.0828"
    $dataQ3[2,6] = 9
    $dataQ3[3,0] = "'003721"
    $dataQ3[3,1] = "易方达标普信息科技指数（QDII-LOF）美元A"
    $dataQ3[3,2] = "'4.84"
    $dataQ3[3,3] = "'91.96"
    $dataQ3[3,4] = "'1.66"
    $dataQ3[3,5] = "'0.0803"
    $dataQ3[3,6] = 9
    $dataQ3[4,0] = "'012869"
    $dataQ3[4,1] = "易方达标普信息科技指数（QDII-LOF）美元 C"
    $dataQ3[4,2] = "'0.15"
    $dataQ3[4,3] = "'91.96"
    $dataQ3[4,4] = "'1.66"
    $dataQ3[4,5] = "'0.0025"
    $dataQ3[4,6] = 9
    $dataQ3[5,0] = "'015203"
    $dataQ3[5,1] = "汇添富全球移动互联灵活配置混合（QDII）D"
    $dataQ3[5,2] = "'0.04"
    $dataQ3[5,3] = "'90.88"
    $dataQ3[5,4] = "'2.86"
    $dataQ3[5,5] = "'0.0011"
    $dataQ3[5,6] = 8
    $dataQ3[6,0] = "'015202"
    $dataQ3[6,1] = "汇添富全球移动互联灵活配置混合（QDII）C"
    $dataQ3[6,2] = "'0.01"
    $dataQ3[6,3] = "'90.88"
    $dataQ3[6,4] = "'2.86"
    $dataQ3[6,5] = "'0.0003"
    $dataQ3[6,6] = 8
    $wsQ3.Range("B2:H8").Value = $dataQ3
    $wsQ3.Range("B2:H8").Style = "Normal"

# --- 4. Column A (row index 0..6) ---
$idxQ3 = New-Object 'object[,]' 7,1
for ($i = 0; $i -lt 7; $i++) { $idxQ3[$i,0] = $i }
$wsQ3.Range("A2:A8").Value = $idxQ3

# --- 5. Re-use the existing bold/bordered header style (style index 2) for
#        the header row and column-A cells, by copy/paste-format from the
#        equivalent, already-styled cells on the "总计" sheet. This avoids
#        creating duplicate style entries. ---
$wsTotal.Range("B1").Copy() | Out-Null
$wsQ3.Range("B1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsTotal.Range("A2").Copy() | Out-Null
$wsQ3.Range("A2:A8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false


# ---------------------------------------------------------------------------
# 6. "总计" sheet: push the existing 2022-Q2 summary row down to row 3, and
#    overwrite row 2 with the new 2022-Q3 summary figures.
# ---------------------------------------------------------------------------

# 6a. Write the former row-2 data (2022-Q2 / 4 / 0.19) into row 3.
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 4
$wsTotal.Range("D3").Value = 0.19
$wsTotal.Range("A3").Value = 1

# match A2's style (bold/border, style index 2) on the new A3 cell
$wsTotal.Range("A2").Copy() | Out-Null
$wsTotal.Range("A3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# 6b. Overwrite row 2 with the new 2022-Q3 summary figures.
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 7
$wsTotal.Range("D2").Value = 0.59

Write-Host "2022-Q3 sheet + 总计 summary row added"
